$wb = $excel.ActiveWorkbook

# --- doors sheet: append new row (d003 / 280) ---
$doors = $wb.Worksheets.Item("doors")
$doors.Range("A4").Value = "d003"
$doors.Range("B4").Value = 280

# --- plywood sheet: update existing stock value for p002 ---
$plywood = $wb.Worksheets.Item("plywood")
$plywood.Range("B3").Value = -10

# --- add new "hardware" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$hardware = $wb.Worksheets.Add($null, $lastSheet)
$hardware.Name = "hardware"

$hardware.Range("A1").Value = "id"
$hardware.Range("B1").Value = "stock"
$hardware.Range("A2").Value = "Hardware 1mm"
$hardware.Range("B2").Value = 180
